$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a brand new row before the old row 7 ("Creazione Spot").
#    This shifts the old rows 7-13 down to 8-14 and duplicates the
#    formatting pattern of row 6 (the row immediately above) into the
#    freshly inserted row 7.
# ------------------------------------------------------------------
$ws.Rows("7:7").Insert()

# Re-apply the exact format pattern of (the still untouched) row 6 onto
# the new row 7, since Insert() already approximates it - this makes
# sure borders/fills/alignment match row 6 exactly.
$ws.Range("B6:F6").Copy()
$ws.Range("B7:F7").PasteSpecial(-4122)

# New activity entry in row 7.
$ws.Range("D7").Value = "Creazione JSON utente"

# ------------------------------------------------------------------
# 2) Turn row 6 ("Gestione Routing") into the highlighted row: give it
#    a yellow fill, a new start date, and tweak the note in F6.
#    NumberFormat is applied LAST on B6 so the engine keeps reusing
#    the builtin date format (id 14) instead of minting a custom one.
# ------------------------------------------------------------------
$ws.Range("B6").Interior.Color = 65535
$ws.Range("B6").Value2 = 42732
$ws.Range("B6").NumberFormat = "mm-dd-yy"

$ws.Range("C6").Interior.Color = 65535
$ws.Range("D6").Interior.Color = 65535

$ws.Range("F6").Interior.Color = 65535
$ws.Range("F6").Value = "Routing - passaggio valore a componenti figli - da pulire"

$ws.Range("E6").Interior.Color = 65535

# ------------------------------------------------------------------
# 3) Keep the selection/active cell consistent with the new layout.
# ------------------------------------------------------------------
$ws.Range("F6").Select()
